$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Increased"
$ws.Range("B3").Value = "apple"
$ws.Range("C3").Value = "\31 52170-case-633"
$ws.Range("D3").Value = "\31 52171-case-640"
$ws.Range("E3").Value = "128 GB"
$ws.Range("F3").Value = "Green"
$ws.Range("G3").Value = " "
$ws.Range("H3").Value = "'1"
$ws.Range("I3").Value = "apple-iphone-15-new"
